# Updates the cryptos price table with the latest snapshot values.
# Numeric-looking text (e.g. "306.61") is written with a leading apostrophe
# so Excel keeps it as text (matching the original inlineStr cells) instead
# of silently converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.446.31'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '2.246.80'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '''306.61'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = '''93.53'
$ws.Range("E6").Value = '  -5.61%  '
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("D8").Value = '''1.01'
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = '''0.523'
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("D10").Value = '''34.62'
$ws.Range("E10").Value = '  -3.15%  '
$ws.Range("D11").Value = '''0.0810'
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").Value = '''7.14'
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.244.50'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.837'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").Value = '''13.58'
$ws.Range("E16").Value = '  -2.24%  '
$ws.Range("D17").Value = '44.142.34'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '''12.34'
$ws.Range("E19").Value = '  -4.23%  '
$ws.Range("D20").Value = '''6.34'
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '''65.64'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '''237.95'
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").Value = '''1.98'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '''38.74'
$ws.Range("E26").Value = '  +4.27%  '
$ws.Range("E27").Value = '  +2.93%  '
$ws.Range("E28").Value = '  -3.29%  '
$ws.Range("D29").Value = '''5.93'
$ws.Range("E29").Value = '  -3.66%  '
$ws.Range("D30").Value = '''20.06'
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").Value = '''152.57'
$ws.Range("E31").Value = '  -2.98%  '
$ws.Range("D32").Value = '''0.0799'
$ws.Range("E32").Value = '  -3.63%  '
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").Value = '''3.09'
$ws.Range("E34").Value = '  -13.49%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '''0.119'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.108'
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '''1.81'
$ws.Range("E37").Value = '  -2.84%  '
$ws.Range("D38").Value = '''3.46'
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("D39").Value = '''14.67'
$ws.Range("E39").Value = '  -4.52%  '
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("D41").Value = '''0.0302'
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D43").Value = '1.733.61'
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("D44").Value = '''80.63'
$ws.Range("E44").Value = '  -9.53%  '
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").Value = '''99.48'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '''4.93'
$ws.Range("E47").Value = '  -4.45%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '''1.61'
$ws.Range("E48").Value = '  +4.00%  '
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '''69.59'
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''55.31'
$ws.Range("E51").Value = '  -0.52%  '
